# Fix: "Error in Increment month application is rectified"
#
# The sheet tracks a monthly arrear report. The pay increment that occurs
# each year had erroneously been applied starting in February instead of
# January. This script shifts the increment back to January for each of
# the affected rows (and fixes the follow-on Jan/Feb-2026 values that
# depended on the wrong timing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => column/value pairs that need correcting
$fixes = @{
    14 = @{ C = 78300; D = 81200; F = 88479; G = 91756; H = 3277 }
    26 = @{ C = 80700; D = 83700; F = 91191; G = 94581; H = 3390 }
    38 = @{ C = 83200; D = 86300; F = 94016; G = 97519; H = 3503 }
    50 = @{ C = 85700; D = 88900; F = 102840; G = 106680; H = 3840 }
    62 = @{ C = 88300; D = 91600; F = 109492; G = 113584; H = 4092 }
    74 = @{ B = 6600; C = 91000; D = 94400; F = 116480; G = 120832; H = 4352 }
    75 = @{ C = 96800; D = 102600; F = 123904; G = 131328; H = 7424 }
}

foreach ($row in $fixes.Keys) {
    $cols = $fixes[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
